$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A24").Value = "edit1"
$ws.Range("B24").Value = "riya-morankar"
$ws.Range("C24").Value = "Merged"
$ws.Range("D24").Value = "N/A"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2025-06-23"
$ws.Range("F24").Value = "7cb3f674f05f8033cb60e20bc0d35a1e90579a4e"
